# feat(commands): adapt xlspython to the new MVC architecture
#
# The "jouet.xlsx" test sheet used to expose 5 columns:
#   A: Prenom, B: Etat, C: Temps utilise, D: Note/10,00, E: Reponse 1
# The new MVC command layer only needs the timing + grade columns, so the
# "Prenom"/"Etat" columns are dropped and the trailing "Reponse 1" column
# (free-text consent answer) is dropped as well, leaving just:
#   A: Temps utilise, B: Note/10,00

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the leading "Prenom" / "Etat" columns - this shifts
# "Temps utilise" -> A, "Note/10,00" -> B, "Reponse 1" -> C.
$ws.Range("A:B").EntireColumn.Delete()

# Drop the trailing "Reponse 1" column (now column C).
$ws.Range("C:C").EntireColumn.Delete()
